# Auto-generated Excel COM-interop script
# Rewrites data rows 2-38 (columns A:M) of Sheet1 so the workbook matches
# the target edit described in the diff (corrected/added FOP rows, region
# names lowercased/abbreviated, surname+region concatenation in col M, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing date-style number format (column G / J use it)
# so newly-written rows get the same style index as the pre-existing ones.
$dateFormat = $ws.Range("G2").NumberFormat

# Columns A (код ЄДРПОУ) and C (Серія препарату) hold text values that can
# look numeric (including values with significant leading zeros, e.g.
# "0401622"). Flip them to Text format before writing, then back to General
# afterwards, so COM stores the literal text instead of silently coercing
# it to a number (which would also strip leading zeros).
$ws.Range("A2:A38").NumberFormat = "@"
$ws.Range("C2:C38").NumberFormat = "@"

# Row 2
$row2 = New-Object 'object[,]' 1,13
$row2[0,0] = '3255920623'
$row2[0,1] = 'ФОП Григоришина'
$row2[0,2] = '1141M019A'
$row2[0,3] = 'ХІБ'
$row2[0,4] = 'ВАКЦИНА КОН`ЮГОВАНА ДЛЯ ПРОФІЛАКТИКИ ЗАХВОРЮВАНЬ ЗБУДНИКОМ ЯКИХ Є HAEMOPHILUS INFLUENZAE ТИПУ В. Ліофілізат для розчину для ін`єкцій по 1 дозі (10 мкг PRP) у флаконах № 50 у комплекті з розчинником ( 04 % розчин натрію хлориду) у ампулах №50,'
$row2[0,5] = 'СЕРУМ ІНСТИТУТ ІНДІЇ ПВТ. ЛТД., Індія'
$row2[0,6] = 45433
$row2[0,7] = 'Державний бюджет 2021'
$row2[0,8] = 11
$row2[0,9] = 45345
$row2[0,10] = 'вінницька'
$row2[0,11] = 'ФОП'
$row2[0,12] = 'григоришина+вінницька'
$ws.Range("A2:M2").Value = $row2
$ws.Range("G2").NumberFormat = $dateFormat
$ws.Range("J2").NumberFormat = $dateFormat

# Row 3
$row3 = New-Object 'object[,]' 1,13
$row3[0,0] = '2335110866'
$row3[0,1] = 'ФОП Червякова'
$row3[0,2] = 'UFA22005'
$row3[0,3] = 'ГепВ'
$row3[0,4] = 'ЕУВАКС В Вакцина для профілактики гепатиту В, рекомбінантна рідка, суспензія для ін''єкцій, 10 мкг/дозу,по 0,5 мл (1 доза) у флаконі; по 10 флаконів у картонній пачці'
$row3[0,5] = 'ЕлДжі Лайф Сайенсіс Лтд,Корея'
$row3[0,6] = 45545
$row3[0,7] = 'Державний бюджет 2021'
$row3[0,8] = 5
$row3[0,9] = 45345
$row3[0,10] = 'вінницька'
$row3[0,11] = 'ФОП'
$row3[0,12] = 'червякова+вінницька'
$ws.Range("A3:M3").Value = $row3
$ws.Range("G3").NumberFormat = $dateFormat
$ws.Range("J3").NumberFormat = $dateFormat

# Row 4
$row4 = New-Object 'object[,]' 1,13
$row4[0,0] = '2673513087'
$row4[0,1] = 'ФОП Г.Кліщевська'
$row4[0,2] = 'X007788'
$row4[0,3] = 'КПК'
$row4[0,4] = 'MMRVAXPRO II 0,5ML 10X1DOSE VIAL L25'
$row4[0,5] = 'Merck & Co., Inc.'
$row4[0,6] = 45432
$row4[0,7] = 'Гуманітарна допомога'
$row4[0,8] = 34
$row4[0,9] = 45345
$row4[0,10] = 'волинська'
$row4[0,11] = 'ФОП'
$row4[0,12] = 'г.кліщевська+волинська'
$ws.Range("A4:M4").Value = $row4
$ws.Range("G4").NumberFormat = $dateFormat
$ws.Range("J4").NumberFormat = $dateFormat

# Row 5
$row5 = New-Object 'object[,]' 1,13
$row5[0,0] = '3166913315'
$row5[0,1] = 'ФОП Михалевський Б.'
$row5[0,2] = '0401622'
$row5[0,3] = 'АДП'
$row5[0,4] = 'АДСОРБОВАНА ВАКЦИНА ДТ / ADSORBED DT VACCINE, суспензія для ін''єкцій по 10 доз (1 доза 0.5 мл) по 5 мл у флаконах № 10 у картонній коробці'
$row5[0,5] = 'ПТ БІО ФАРМА (Персеро), Індонезія'
$row5[0,6] = 46192
$row5[0,7] = 'Державний бюджет 2021'
$row5[0,8] = 19
$row5[0,9] = 45345
$row5[0,10] = 'волинська'
$row5[0,11] = 'ФОП'
$row5[0,12] = 'михалевський+волинська'
$ws.Range("A5:M5").Value = $row5
$ws.Range("G5").NumberFormat = $dateFormat
$ws.Range("J5").NumberFormat = $dateFormat

# Row 6
$row6 = New-Object 'object[,]' 1,13
$row6[0,0] = '2733508291'
$row6[0,1] = 'ФОП "Горошко А.М."'
$row6[0,2] = 'D2720-E1'
$row6[0,3] = 'АДП-М'
$row6[0,4] = 'ТЕТАДІФ суспензія для ін''єкцій 0,5 мл (1 доза), флакони по 5 мл. (10 доз), по 10 флаконів у картонній коробці'
$row6[0,5] = 'ББ-НЦІПХЛтд. Болгарія'
$row6[0,6] = 46081
$row6[0,7] = 'Державний бюджет 2021'
$row6[0,8] = 46
$row6[0,9] = 45345
$row6[0,10] = 'дніпропетровська'
$row6[0,11] = 'ФОП'
$row6[0,12] = '"горошко+дніпропетровська'
$ws.Range("A6:M6").Value = $row6
$ws.Range("G6").NumberFormat = $dateFormat
$ws.Range("J6").NumberFormat = $dateFormat

# Row 7
$row7 = New-Object 'object[,]' 1,13
$row7[0,0] = '3254011678'
$row7[0,1] = 'ФОП Білишко'
$row7[0,2] = '0401622'
$row7[0,3] = 'АДП'
$row7[0,4] = 'АДСОРБОВАНА ВАКЦИНА ДТ / ADSORBED DT VACCINE, суспензія для ін''єкцій по 10 доз (1 доза 0.5 мл) по 5 мл у флаконах № 10 у картонній коробці'
$row7[0,5] = 'ПТ БІО ФАРМА (Персеро), Індонезія'
$row7[0,6] = 46012
$row7[0,7] = 'Державний бюджет 2021'
$row7[0,8] = 100
$row7[0,9] = 45345
$row7[0,10] = 'дніпропетровська'
$row7[0,11] = 'ФОП'
$row7[0,12] = 'білишко+дніпропетровська'
$ws.Range("A7:M7").Value = $row7
$ws.Range("G7").NumberFormat = $dateFormat
$ws.Range("J7").NumberFormat = $dateFormat

# Row 8
$row8 = New-Object 'object[,]' 1,13
$row8[0,0] = '2889703358'
$row8[0,1] = 'ФОП "Шпорт А.Я."'
$row8[0,2] = 'FWL21506'
$row8[0,3] = 'Пента'
$row8[0,4] = 'Eupenta inj.DTP-Hib-HepB 0,5mL'
$row8[0,5] = 'ЕлДжі Лайф Сайенсіс Лтд,Корея'
$row8[0,6] = 45657
$row8[0,7] = 'Гуманітарна допомога'
$row8[0,8] = 48
$row8[0,9] = 45345
$row8[0,10] = 'дніпропетровська'
$row8[0,11] = 'ФОП'
$row8[0,12] = '"шпорт+дніпропетровська'
$ws.Range("A8:M8").Value = $row8
$ws.Range("G8").NumberFormat = $dateFormat
$ws.Range("J8").NumberFormat = $dateFormat

# Row 9
$row9 = New-Object 'object[,]' 1,13
$row9[0,0] = '3216713340'
$row9[0,1] = 'ФОП Мерзляк'
$row9[0,2] = '0401622'
$row9[0,3] = 'АДП'
$row9[0,4] = 'АДСОРБОВАНА ВАКЦИНА ДТ / ADSORBED DT VACCINE, суспензія для ін''єкцій по 10 доз (1 доза 0.5 мл) по 5 мл у флаконах № 10 у картонній коробці'
$row9[0,5] = 'ПТ БІО ФАРМА (Персеро), Індонезія'
$row9[0,6] = 46012
$row9[0,7] = 'Державний бюджет 2021'
$row9[0,8] = 60
$row9[0,9] = 45345
$row9[0,10] = 'дніпропетровська'
$row9[0,11] = 'ФОП'
$row9[0,12] = 'мерзляк+дніпропетровська'
$ws.Range("A9:M9").Value = $row9
$ws.Range("G9").NumberFormat = $dateFormat
$ws.Range("J9").NumberFormat = $dateFormat

# Row 10
$row10 = New-Object 'object[,]' 1,13
$row10[0,0] = '2986403644'
$row10[0,1] = 'ФОП Климук К.В.'
$row10[0,2] = '68D23040A'
$row10[0,3] = 'ОПВ'
$row10[0,4] = 'Bivalent Poliomyelitis Vaccine Tipe 1&3, Live (Oral), BIOPOLIO B1/3, 1mL Vail of 10 doses'
$row10[0,5] = 'Bharat Biotech Int.Ltd, Індія'
$row10[0,6] = 45505
$row10[0,7] = 'Гуманітарна допомога'
$row10[0,8] = 10
$row10[0,9] = 45345
$row10[0,10] = 'дніпропетровська'
$row10[0,11] = 'ФОП'
$row10[0,12] = 'климук+дніпропетровська'
$ws.Range("A10:M10").Value = $row10
$ws.Range("G10").NumberFormat = $dateFormat
$ws.Range("J10").NumberFormat = $dateFormat

# Row 11
$row11 = New-Object 'object[,]' 1,13
$row11[0,0] = $null
$row11[0,1] = 'ФОП Пасько Г.І.'
$row11[0,2] = '68D23040А'
$row11[0,3] = 'ОПВ'
$row11[0,4] = 'Bivalent Poliomyelitis Vaccine Tipe 1&3, Live (Oral), BIOPOLIO B1/3, 1mL Vail of 10 doses'
$row11[0,5] = 'Bharat Biotech Int.Ltd, Індія'
$row11[0,6] = 45505
$row11[0,7] = 'Гуманітарна допомога'
$row11[0,8] = 30
$row11[0,9] = 45345
$row11[0,10] = 'дніпропетровська'
$row11[0,11] = 'ФОП'
$row11[0,12] = 'пасько+дніпропетровська'
$ws.Range("A11:M11").Value = $row11
$ws.Range("G11").NumberFormat = $dateFormat
$ws.Range("J11").NumberFormat = $dateFormat

# Row 12
$row12 = New-Object 'object[,]' 1,13
$row12[0,0] = '2562701239'
$row12[0,1] = 'ФОП "Маланчик"'
$row12[0,2] = '2249001B'
$row12[0,3] = 'ІПВ'
$row12[0,4] = 'Inactivated Poliomyelitis Vaccine'
$row12[0,5] = 'Bilthoven Biologicals B.V.'
$row12[0,6] = 45999
$row12[0,7] = 'Гуманітарна допомога'
$row12[0,8] = 11
$row12[0,9] = 45338
$row12[0,10] = 'донецька'
$row12[0,11] = 'ФОП'
$row12[0,12] = '"маланчик"+донецька'
$ws.Range("A12:M12").Value = $row12
$ws.Range("G12").NumberFormat = $dateFormat
$ws.Range("J12").NumberFormat = $dateFormat

# Row 13
$row13 = New-Object 'object[,]' 1,13
$row13[0,0] = '2584015708'
$row13[0,1] = 'ФОП Венгер С.І '
$row13[0,2] = '2204008A'
$row13[0,3] = 'ІПВ'
$row13[0,4] = 'Inactivated poliomyelitis vaccine'
$row13[0,5] = 'Bilthoven Biologicals B.V.A. Netherlands'
$row13[0,6] = 45646
$row13[0,7] = 'Гуманітарна допомога'
$row13[0,8] = 20
$row13[0,9] = 45345
$row13[0,10] = 'житомирська'
$row13[0,11] = 'ФОП'
$row13[0,12] = 'венгер+житомирська'
$ws.Range("A13:M13").Value = $row13
$ws.Range("G13").NumberFormat = $dateFormat
$ws.Range("J13").NumberFormat = $dateFormat

# Row 14
$row14 = New-Object 'object[,]' 1,13
$row14[0,0] = '3451917236'
$row14[0,1] = 'ФОП Кабан В.Б.'
$row14[0,2] = '2310002A'
$row14[0,3] = 'ІПВ'
$row14[0,4] = 'Inactivated Poliomyelitis Vaccine'
$row14[0,5] = 'Inactivated Poliomyelitis Vaccine'
$row14[0,6] = 46067
$row14[0,7] = 'Гуманітарна допомога'
$row14[0,8] = 30
$row14[0,9] = 45331
$row14[0,10] = 'івано-франківська'
$row14[0,11] = 'ФОП'
$row14[0,12] = 'кабан+івано-франківська'
$ws.Range("A14:M14").Value = $row14
$ws.Range("G14").NumberFormat = $dateFormat
$ws.Range("J14").NumberFormat = $dateFormat

# Row 15
$row15 = New-Object 'object[,]' 1,13
$row15[0,0] = '43968084'
$row15[0,1] = 'ФОП "Бенько Г.С."'
$row15[0,2] = 'D2708'
$row15[0,3] = 'АДП-М'
$row15[0,4] = 'ТЕТАДІФ суспензія для ін''єкцій 0,5 мл (1 доза), флакони по 5 мл. (10 доз), по 10 флаконів у картонній коробці'
$row15[0,5] = 'ТЕТАДІФ суспензія для ін''єкцій 0,5 мл (1 доза), флакони по 5 мл. (10 доз), по 10 флаконів у картонній коробці'
$row15[0,6] = 45991
$row15[0,7] = 'Державний бюджет 2021'
$row15[0,8] = 17
$row15[0,9] = 45331
$row15[0,10] = 'івано-франківська'
$row15[0,11] = 'ФОП'
$row15[0,12] = '"бенько+івано-франківська'
$ws.Range("A15:M15").Value = $row15
$ws.Range("G15").NumberFormat = $dateFormat
$ws.Range("J15").NumberFormat = $dateFormat

# Row 16
$row16 = New-Object 'object[,]' 1,13
$row16[0,0] = '3334619475'
$row16[0,1] = 'ФОП Венгринович В.В.'
$row16[0,2] = 'D2708'
$row16[0,3] = 'АДП-М'
$row16[0,4] = 'ТЕТАДІФ суспензія для ін''єкцій 0,5 мл (1 доза), флакони по 5 мл. (10 доз), по 10 флаконів у картонній коробці'
$row16[0,5] = 'ТЕТАДІФ суспензія для ін''єкцій 0,5 мл (1 доза), флакони по 5 мл. (10 доз), по 10 флаконів у картонній коробці'
$row16[0,6] = 45991
$row16[0,7] = 'Державний бюджет 2021'
$row16[0,8] = 30
$row16[0,9] = 45331
$row16[0,10] = 'івано-франківська'
$row16[0,11] = 'ФОП'
$row16[0,12] = 'венгринович+івано-франківська'
$ws.Range("A16:M16").Value = $row16
$ws.Range("G16").NumberFormat = $dateFormat
$ws.Range("J16").NumberFormat = $dateFormat

# Row 17
$row17 = New-Object 'object[,]' 1,13
$row17[0,0] = '2339913111'
$row17[0,1] = 'ФОП "Сем''янчук В.Б."'
$row17[0,2] = 'UFA23503'
$row17[0,3] = 'ГепВ'
$row17[0,4] = 'ЕУВАКС В ВАКЦИНА ДЛЯ ПРОФІЛАКТИКИ ГЕПАТИТУ В РЕКОМБІНАНТНА РІДКА суспензія для ін''єкцій 10 мкг/дозу по 05 мл (1 доза) у флаконі,'
$row17[0,5] = 'ЕУВАКС В ВАКЦИНА ДЛЯ ПРОФІЛАКТИКИ ГЕПАТИТУ В РЕКОМБІНАНТНА РІДКА суспензія для ін''єкцій 10 мкг/дозу по 05 мл (1 доза) у флаконі,'
$row17[0,6] = 46088
$row17[0,7] = 'Державний бюджет 2022'
$row17[0,8] = 20
$row17[0,9] = 45331
$row17[0,10] = 'івано-франківська'
$row17[0,11] = 'ФОП'
$row17[0,12] = '"сем''янчук+івано-франківська'
$ws.Range("A17:M17").Value = $row17
$ws.Range("G17").NumberFormat = $dateFormat
$ws.Range("J17").NumberFormat = $dateFormat

# Row 18
$row18 = New-Object 'object[,]' 1,13
$row18[0,0] = '2859407383'
$row18[0,1] = 'ФОП Мегединник Ганна Володимирівна'
$row18[0,2] = 'X007787'
$row18[0,3] = 'КПК'
$row18[0,4] = 'MMRVAXPRO II 0,5ML 10X1DOSE VIAL L25'
$row18[0,5] = 'MMRVAXPRO II 0,5ML 10X1DOSE VIAL L25'
$row18[0,6] = 45438
$row18[0,7] = 'Гуманітарна допомога'
$row18[0,8] = 113
$row18[0,9] = 45331
$row18[0,10] = 'івано-франківська'
$row18[0,11] = 'ФОП'
$row18[0,12] = 'мегединник+івано-франківська'
$ws.Range("A18:M18").Value = $row18
$ws.Range("G18").NumberFormat = $dateFormat
$ws.Range("J18").NumberFormat = $dateFormat

# Row 19
$row19 = New-Object 'object[,]' 1,13
$row19[0,0] = '3076817858'
$row19[0,1] = 'ФОП Білошицький Олександр Михайлович'
$row19[0,2] = 'AMJRE616AB'
$row19[0,3] = 'КПК'
$row19[0,4] = 'PRIO.VIAL.2DX100 VVM DX WH (TEND) Combined Measles, Mumps and Rubella vaccine (live) '
$row19[0,5] = 'PRIO.VIAL.2DX100 VVM DX WH (TEND) Combined Measles, Mumps and Rubella vaccine (live) '
$row19[0,6] = 45473
$row19[0,7] = 'Державний бюджет 2021'
$row19[0,8] = 14
$row19[0,9] = 45331
$row19[0,10] = 'івано-франківська'
$row19[0,11] = 'ФОП'
$row19[0,12] = 'білошицький+івано-франківська'
$ws.Range("A19:M19").Value = $row19
$ws.Range("G19").NumberFormat = $dateFormat
$ws.Range("J19").NumberFormat = $dateFormat

# Row 20
$row20 = New-Object 'object[,]' 1,13
$row20[0,0] = '2339913117'
$row20[0,1] = 'Здорова малеча ФОП Савюк '
$row20[0,2] = '2310002A'
$row20[0,3] = 'ІПВ'
$row20[0,4] = 'Inactivated Poliomyelitis Vaccine'
$row20[0,5] = 'Inactivated Poliomyelitis Vaccine'
$row20[0,6] = 46067
$row20[0,7] = 'Гуманітарна допомога'
$row20[0,8] = 22
$row20[0,9] = 45331
$row20[0,10] = 'івано-франківська'
$row20[0,11] = 'ФОП'
$row20[0,12] = 'здорова+івано-франківська'
$ws.Range("A20:M20").Value = $row20
$ws.Range("G20").NumberFormat = $dateFormat
$ws.Range("J20").NumberFormat = $dateFormat

# Row 21
$row21 = New-Object 'object[,]' 1,13
$row21[0,0] = '3285120231'
$row21[0,1] = 'ФОП Корнеев'
$row21[0,2] = 'AMJRE680AA'
$row21[0,3] = 'КПК'
$row21[0,4] = 'PRIORIX vaccine'
$row21[0,5] = 'Гласко Сміт Кляйн,Біолоджікалз С.А.,Бельгія'
$row21[0,6] = 45657
$row21[0,7] = 'Гуманітарна допомога'
$row21[0,8] = 2
$row21[0,9] = 45345
$row21[0,10] = 'київська'
$row21[0,11] = 'ФОП'
$row21[0,12] = 'корнеев+київська'
$ws.Range("A21:M21").Value = $row21
$ws.Range("G21").NumberFormat = $dateFormat
$ws.Range("J21").NumberFormat = $dateFormat

# Row 22
$row22 = New-Object 'object[,]' 1,13
$row22[0,0] = '2936213862'
$row22[0,1] = 'ФОП Мороз'
$row22[0,2] = 'D2720-E1'
$row22[0,3] = 'АДП-М'
$row22[0,4] = 'ТЕТАДІФ суспензія для ін''єкцій 0,5 мл (1 доза), флакони по 5 мл. (10 доз), по 10 флаконів у картонній коробці'
$row22[0,5] = 'ББ-НЦІПХЛтд. Болгарія'
$row22[0,6] = 46081
$row22[0,7] = 'Державний бюджет 2021'
$row22[0,8] = 10
$row22[0,9] = 45345
$row22[0,10] = 'київська'
$row22[0,11] = 'ФОП'
$row22[0,12] = 'мороз+київська'
$ws.Range("A22:M22").Value = $row22
$ws.Range("G22").NumberFormat = $dateFormat
$ws.Range("J22").NumberFormat = $dateFormat

# Row 23
$row23 = New-Object 'object[,]' 1,13
$row23[0,0] = '3198123020'
$row23[0,1] = 'ФОП Лавренчук Ірина Олегівна'
$row23[0,2] = '68D23040А'
$row23[0,3] = 'ОПВ'
$row23[0,4] = 'Bivalent Poliomyelitis Vaccine Tipe 1&3, Live (Oral), BIOPOLIO B1/3, 1mL Vail of 10 doses'
$row23[0,5] = 'Bharat Biotech Int.Ltd, Індія'
$row23[0,6] = 45505
$row23[0,7] = 'Гуманітарна допомога'
$row23[0,8] = 47
$row23[0,9] = 45345
$row23[0,10] = 'київська'
$row23[0,11] = 'ФОП'
$row23[0,12] = 'лавренчук+київська'
$ws.Range("A23:M23").Value = $row23
$ws.Range("G23").NumberFormat = $dateFormat
$ws.Range("J23").NumberFormat = $dateFormat

# Row 24
$row24 = New-Object 'object[,]' 1,13
$row24[0,0] = '3311003792'
$row24[0,1] = 'ФОП ПАГІЄВ І.Ф.'
$row24[0,2] = '0401622'
$row24[0,3] = 'АДП'
$row24[0,4] = 'АДСОРБОВАНА ВАКЦИНА ДТ / ADSORBED DT VACCINE, суспензія для ін''єкцій по 10 доз (1 доза 0.5 мл) по 5 мл у флаконах № 10 у картонній коробці'
$row24[0,5] = 'ПТ БІО ФАРМА (Персеро), Індонезія'
$row24[0,6] = 46012
$row24[0,7] = 'Державний бюджет 2021'
$row24[0,8] = 38
$row24[0,9] = 45345
$row24[0,10] = 'кіровоградська'
$row24[0,11] = 'ФОП'
$row24[0,12] = 'пагієв+кіровоградська'
$ws.Range("A24:M24").Value = $row24
$ws.Range("G24").NumberFormat = $dateFormat
$ws.Range("J24").NumberFormat = $dateFormat

# Row 25
$row25 = New-Object 'object[,]' 1,13
$row25[0,0] = '3312303896'
$row25[0,1] = 'ФОП ЯКОВЕНКО'
$row25[0,2] = '222600722B'
$row25[0,3] = 'АДП-М'
$row25[0,4] = 'ВАКЦИНА ДЛЯ ПРОФІЛАКТИКИ ДИФТЕРІЇ ТА ПРАВЦЯ АДСОРБОВАНА ІЗ ЗМЕНШЕНИМ ВМІСТОМ АНТИГЕНУ суспензія для ін''єкцій по 10 доз (одна доза 05 мл) по 5 мл у флаконі; по 24 флакони в пачці з картону,'
$row25[0,5] = 'БАЙОЛОДЖІКАЛ І. ЛІМІТЕД, Індія'
$row25[0,6] = 45777
$row25[0,7] = 'Гуманітарна допомога'
$row25[0,8] = 60
$row25[0,9] = 45345
$row25[0,10] = 'кіровоградська'
$row25[0,11] = 'ФОП'
$row25[0,12] = 'яковенко+кіровоградська'
$ws.Range("A25:M25").Value = $row25
$ws.Range("G25").NumberFormat = $dateFormat
$ws.Range("J25").NumberFormat = $dateFormat

# Row 26
$row26 = New-Object 'object[,]' 1,13
$row26[0,0] = '3088718663'
$row26[0,1] = 'ФОП "РАДЬКО"'
$row26[0,2] = '2324002A'
$row26[0,3] = 'ІПВ'
$row26[0,4] = 'Inactivated Poliomyelitis Vaccine'
$row26[0,5] = 'Bilthoven Biologicals B.V. Netherlands'
$row26[0,6] = 46155
$row26[0,7] = 'Гуманітарна допомога'
$row26[0,8] = 5
$row26[0,9] = 45345
$row26[0,10] = 'кіровоградська'
$row26[0,11] = 'ФОП'
$row26[0,12] = '"радько"+кіровоградська'
$ws.Range("A26:M26").Value = $row26
$ws.Range("G26").NumberFormat = $dateFormat
$ws.Range("J26").NumberFormat = $dateFormat

# Row 27
$row27 = New-Object 'object[,]' 1,13
$row27[0,0] = '3296511462'
$row27[0,1] = 'ФОП ПЕНЬКАЧ'
$row27[0,2] = '221100223B'
$row27[0,3] = 'АКДП'
$row27[0,4] = 'DTP VACCINE (2023) Adsorbed Diphtheria, Tetanus and Pertussis vaccine Pediatric dose, vial of 10 doses
With Vaccine Vial Monitor (VVM)'
$row27[0,5] = 'БАЙОЛОДЖІКАЛ І. ЛІМІТЕД, Індія'
$row27[0,6] = 45807
$row27[0,7] = 'Гуманітарна допомога'
$row27[0,8] = 129
$row27[0,9] = 45345
$row27[0,10] = 'кіровоградська'
$row27[0,11] = 'ФОП'
$row27[0,12] = 'пенькач+кіровоградська'
$ws.Range("A27:M27").Value = $row27
$ws.Range("G27").NumberFormat = $dateFormat
$ws.Range("J27").NumberFormat = $dateFormat

# Row 28
$row28 = New-Object 'object[,]' 1,13
$row28[0,0] = '2918018806'
$row28[0,1] = 'ФОП Заремба '
$row28[0,2] = 'X007787'
$row28[0,3] = 'КПК'
$row28[0,4] = 'MMRVAXPRO II 0,5ML 10X1DOSE VIAL L25'
$row28[0,5] = 'Merck & Co., Inc.'
$row28[0,6] = 45438
$row28[0,7] = 'Гуманітарна допомога'
$row28[0,8] = 35
$row28[0,9] = 45338
$row28[0,10] = 'львівська'
$row28[0,11] = 'ФОП'
$row28[0,12] = 'заремба+львівська'
$ws.Range("A28:M28").Value = $row28
$ws.Range("G28").NumberFormat = $dateFormat
$ws.Range("J28").NumberFormat = $dateFormat

# Row 29
$row29 = New-Object 'object[,]' 1,13
$row29[0,0] = '3432003184'
$row29[0,1] = 'ФОП Ходан '
$row29[0,2] = '2204008A'
$row29[0,3] = 'ІПВ'
$row29[0,4] = 'Inactivated poliomyelitis vaccine'
$row29[0,5] = 'Bilthoven Biologicals B.V.A. Netherlands'
$row29[0,6] = 45646
$row29[0,7] = 'Гуманітарна допомога'
$row29[0,8] = 15
$row29[0,9] = 45338
$row29[0,10] = 'львівська'
$row29[0,11] = 'ФОП'
$row29[0,12] = 'ходан+львівська'
$ws.Range("A29:M29").Value = $row29
$ws.Range("G29").NumberFormat = $dateFormat
$ws.Range("J29").NumberFormat = $dateFormat

# Row 30
$row30 = New-Object 'object[,]' 1,13
$row30[0,0] = '2481607082'
$row30[0,1] = 'ФОП "Качуровська Ж.Д."'
$row30[0,2] = '1141M019A'
$row30[0,3] = 'ХІБ'
$row30[0,4] = 'ВАКЦИНА КОН`ЮГОВАНА ДЛЯ ПРОФІЛАКТИКИ ЗАХВОРЮВАНЬ ЗБУДНИКОМ ЯКИХ Є HAEMOPHILUS INFLUENZAE ТИПУ В. Ліофілізат для розчину для ін`єкцій по 1 дозі (10 мкг PRP) у флаконах № 50 у комплекті з розчинником ( 04 % розчин натрію хлориду) у ампулах №50,'
$row30[0,5] = 'СЕРУМ ІНСТИТУТ ІНДІЇ ПВТ. ЛТД., Індія'
$row30[0,6] = 45433
$row30[0,7] = 'Державний бюджет 2021'
$row30[0,8] = 1
$row30[0,9] = 45345
$row30[0,10] = 'миколаївська'
$row30[0,11] = 'ФОП'
$row30[0,12] = '"качуровська+миколаївська'
$ws.Range("A30:M30").Value = $row30
$ws.Range("G30").NumberFormat = $dateFormat
$ws.Range("J30").NumberFormat = $dateFormat

# Row 31
$row31 = New-Object 'object[,]' 1,13
$row31[0,0] = '2847608002'
$row31[0,1] = 'ФОП "Коровкіна К.Ю."'
$row31[0,2] = '1141M019A'
$row31[0,3] = 'ХІБ'
$row31[0,4] = 'ВАКЦИНА КОН`ЮГОВАНА ДЛЯ ПРОФІЛАКТИКИ ЗАХВОРЮВАНЬ ЗБУДНИКОМ ЯКИХ Є HAEMOPHILUS INFLUENZAE ТИПУ В. Ліофілізат для розчину для ін`єкцій по 1 дозі (10 мкг PRP) у флаконах № 50 у комплекті з розчинником ( 04 % розчин натрію хлориду) у ампулах №50,'
$row31[0,5] = 'СЕРУМ ІНСТИТУТ ІНДІЇ ПВТ. ЛТД., Індія'
$row31[0,6] = 45433
$row31[0,7] = 'Державний бюджет 2021'
$row31[0,8] = 1
$row31[0,9] = 45345
$row31[0,10] = 'миколаївська'
$row31[0,11] = 'ФОП'
$row31[0,12] = '"коровкіна+миколаївська'
$ws.Range("A31:M31").Value = $row31
$ws.Range("G31").NumberFormat = $dateFormat
$ws.Range("J31").NumberFormat = $dateFormat

# Row 32
$row32 = New-Object 'object[,]' 1,13
$row32[0,0] = '252818590'
$row32[0,1] = 'ФОП "Тімнов В.О."'
$row32[0,2] = '1141M019A'
$row32[0,3] = 'ХІБ'
$row32[0,4] = 'ВАКЦИНА КОН`ЮГОВАНА ДЛЯ ПРОФІЛАКТИКИ ЗАХВОРЮВАНЬ ЗБУДНИКОМ ЯКИХ Є HAEMOPHILUS INFLUENZAE ТИПУ В. Ліофілізат для розчину для ін`єкцій по 1 дозі (10 мкг PRP) у флаконах № 50 у комплекті з розчинником ( 04 % розчин натрію хлориду) у ампулах №50,'
$row32[0,5] = 'СЕРУМ ІНСТИТУТ ІНДІЇ ПВТ. ЛТД., Індія'
$row32[0,6] = 45433
$row32[0,7] = 'Державний бюджет 2021'
$row32[0,8] = 1
$row32[0,9] = 45345
$row32[0,10] = 'миколаївська'
$row32[0,11] = 'ФОП'
$row32[0,12] = '"тімнов+миколаївська'
$ws.Range("A32:M32").Value = $row32
$ws.Range("G32").NumberFormat = $dateFormat
$ws.Range("J32").NumberFormat = $dateFormat

# Row 33
$row33 = New-Object 'object[,]' 1,13
$row33[0,0] = '3400210471'
$row33[0,1] = 'Мій лікар ФОП Гузик Владислав Олегович'
$row33[0,2] = '222600722B'
$row33[0,3] = 'АДП-М'
$row33[0,4] = 'ВАКЦИНА ДЛЯ ПРОФІЛАКТИКИ ДИФТЕРІЇ ТА ПРАВЦЯ АДСОРБОВАНА ІЗ ЗМЕНШЕНИМ ВМІСТОМ АНТИГЕНУ суспензія для ін''єкцій по 10 доз (одна доза 05 мл) по 5 мл у флаконі; по 24 флакони в пачці з картону,'
$row33[0,5] = 'БАЙОЛОДЖІКАЛ І. ЛІМІТЕД, Індія'
$row33[0,6] = 45777
$row33[0,7] = 'Гуманітарна допомога'
$row33[0,8] = 5
$row33[0,9] = 45338
$row33[0,10] = 'одеська'
$row33[0,11] = 'ФОП'
$row33[0,12] = 'мій+одеська'
$ws.Range("A33:M33").Value = $row33
$ws.Range("G33").NumberFormat = $dateFormat
$ws.Range("J33").NumberFormat = $dateFormat

# Row 34
$row34 = New-Object 'object[,]' 1,13
$row34[0,0] = '2519414605'
$row34[0,1] = 'ФОП ЗПСМ та терапія Ланкіна Г.І.'
$row34[0,2] = 'X008927'
$row34[0,3] = 'КПК'
$row34[0,4] = 'MMRVAXPRO II 0,5ML 10X1DOSE VIAL L25'
$row34[0,5] = 'Merck & Co., Inc.'
$row34[0,6] = 45438
$row34[0,7] = 'Гуманітарна допомога'
$row34[0,8] = 30
$row34[0,9] = 45338
$row34[0,10] = 'одеська'
$row34[0,11] = 'ФОП'
$row34[0,12] = 'зпсм+одеська'
$ws.Range("A34:M34").Value = $row34
$ws.Range("G34").NumberFormat = $dateFormat
$ws.Range("J34").NumberFormat = $dateFormat

# Row 35
$row35 = New-Object 'object[,]' 1,13
$row35[0,0] = '2598520904'
$row35[0,1] = 'ФОП Соколова Лариса Володимирівна'
$row35[0,2] = 'AMJRE609BA'
$row35[0,3] = 'КПК'
$row35[0,4] = 'PRIO.VIAL.2DX100 VVM DX WH (TEND) Combined Measles, Mumps and Rubella vaccine (live)'
$row35[0,5] = 'Гласко Сміт Кляйн,Біолоджікалз С.А.,Бельгія'
$row35[0,6] = 45473
$row35[0,7] = 'Гуманітарна допомога'
$row35[0,8] = 8
$row35[0,9] = 45345
$row35[0,10] = 'полтавська'
$row35[0,11] = 'ФОП'
$row35[0,12] = 'соколова+полтавська'
$ws.Range("A35:M35").Value = $row35
$ws.Range("G35").NumberFormat = $dateFormat
$ws.Range("J35").NumberFormat = $dateFormat

# Row 36
$row36 = New-Object 'object[,]' 1,13
$row36[0,0] = '2727016128'
$row36[0,1] = 'ФОП КовенькоТ.'
$row36[0,2] = 'X007039'
$row36[0,3] = 'КПК'
$row36[0,4] = 'MMRVAXPRO II 0,5ML 10X1DOSE VIAL L25'
$row36[0,5] = 'Merck & Co., Inc.'
$row36[0,6] = 45725
$row36[0,7] = 'Гуманітарна допомога'
$row36[0,8] = 19
$row36[0,9] = 45344
$row36[0,10] = 'рівненська'
$row36[0,11] = 'ФОП'
$row36[0,12] = 'ковенькот.+рівненська'
$ws.Range("A36:M36").Value = $row36
$ws.Range("G36").NumberFormat = $dateFormat
$ws.Range("J36").NumberFormat = $dateFormat

# Row 37
$row37 = New-Object 'object[,]' 1,13
$row37[0,0] = '3302400393'
$row37[0,1] = 'ФОП Млавець С.Ф.'
$row37[0,2] = '68D23040А'
$row37[0,3] = 'ОПВ'
$row37[0,4] = 'Bivalent Poliomyelitis Vaccine Tipe 1&3, Live (Oral), BIOPOLIO B1/3, 1mL Vail of 10 doses'
$row37[0,5] = 'Bharat Biotech Int.Ltd, Індія'
$row37[0,6] = 45505
$row37[0,7] = 'Гуманітарна допомога'
$row37[0,8] = 20
$row37[0,9] = 45344
$row37[0,10] = 'черкаська'
$row37[0,11] = 'ФОП'
$row37[0,12] = 'млавець+черкаська'
$ws.Range("A37:M37").Value = $row37
$ws.Range("G37").NumberFormat = $dateFormat
$ws.Range("J37").NumberFormat = $dateFormat

# Row 38
$row38 = New-Object 'object[,]' 1,13
$row38[0,0] = '3286713595'
$row38[0,1] = 'ФОП Агафонов Є.О.'
$row38[0,2] = '68D23011A'
$row38[0,3] = 'ОПВ'
$row38[0,4] = 'Bivalent Poliomyelitis Vaccine Tipe 1&3, Live (Oral), BIOPOLIO B1/3, 1mL Vail of 10 doses'
$row38[0,5] = 'Bharat Biotech Int.Ltd, Індія'
$row38[0,6] = 45377
$row38[0,7] = 'Гуманітарна допомога'
$row38[0,8] = 17
$row38[0,9] = 45345
$row38[0,10] = 'запорізька'
$row38[0,11] = 'ФОП'
$row38[0,12] = 'агафонов+запорізька'
$ws.Range("A38:M38").Value = $row38
$ws.Range("G38").NumberFormat = $dateFormat
$ws.Range("J38").NumberFormat = $dateFormat

# Restore General format on columns A and C now that the text values are in place.
$ws.Range("A2:A38").NumberFormat = "General"
$ws.Range("C2:C38").NumberFormat = "General"

